$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Seshamalini Mohan"
$ws.Range("A5").Value = "Lekha Tummala"

$ws.Range("A6").Select()
